$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "referece" -> "reference" typo in the 10 "libelle_valeur_reference_N" headers.
$ws.Range("B1").Value = "libelle_valeur_reference_1"
$ws.Range("E1").Value = "libelle_valeur_reference_2"
$ws.Range("H1").Value = "libelle_valeur_reference_3"
$ws.Range("K1").Value = "libelle_valeur_reference_4"
$ws.Range("N1").Value = "libelle_valeur_reference_5"
$ws.Range("Q1").Value = "libelle_valeur_reference_6"
$ws.Range("T1").Value = "libelle_valeur_reference_7"
$ws.Range("W1").Value = "libelle_valeur_reference_8"
$ws.Range("Z1").Value = "libelle_valeur_reference_9"
$ws.Range("AC1").Value = "libelle_valeur_reference_10"

# Move the view back to the top-left of the sheet and select C10.
$ws.Range("A1").Select()
$ws.Range("C10").Select()
